$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "ruta"/"rutaDocs" paths (Desktop -> Documents) ---
$ws.Range("B2").Value = "C:/Documents/Tita/"
$ws.Range("B3").Value = "C:/Documents/Tita/Docs/"

# --- Update the contact e-mails (rows 8-11) and mirror them into column C ---
# Capture the hyperlink cell style before the hyperlinks are removed so the
# new column C cells can be styled the same way.
$hlStyle = $ws.Range("B8").Style

$ws.Range("B8").Value = "rafael@2nv.co"
$ws.Range("C8").Value = "rafael@2nv.co"
$ws.Range("C8").Style = $hlStyle

$ws.Range("B9").Value = "jm@2nv.co"
$ws.Range("C9").Value = "jm@2nv.co"
$ws.Range("C9").Style = $hlStyle

$ws.Range("B10").Value = "saul@2nv.co"
$ws.Range("C10").Value = "saul@2nv.co"
$ws.Range("C10").Style = $hlStyle

$ws.Range("B11").Value = "saul@2nv.co"
$ws.Range("C11").Value = "saul@2nv.co"
$ws.Range("C11").Style = $hlStyle

# --- Drop the mailto hyperlinks (values are now plain text) ---
$ws.Hyperlinks.Delete()

# --- Update the active selection shown when the sheet is opened ---
$ws.Range("C8").Select()
